$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.895.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.57"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.60"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0903"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.621.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.561"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.916.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.52"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0700"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.78"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.43"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.109"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.55"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0488"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.16"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.423.81"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.63%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.31"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +8.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.555"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0497"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.826"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "50.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.773.05"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.31"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0112"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.94%  "
